$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header typo: "locacalizacion" -> "localizacion"
$ws.Range("B1").Value = "localizacion"

# Row 3 now holds "manuel"'s data instead of a duplicate of "jorge"'s data
$ws.Range("A3").Value = "manuel"
$ws.Range("B3").Value = "18:15:14:12S"
$ws.Range("C3").Value = "manuel@email.es"
$ws.Range("D3").Value = "ID5"
$ws.Range("E3").Value = 2

# C3 gets a hyperlink to the e-mail address it displays
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:manuel@email.es")
$ws.Range("C3").Style = "Hipervínculo"

# Move the selection to E3
$ws.Range("E3").Select()
